# TC02_Trials_Filter_AssocFileType-AlignedRNA.xlsx
# "commiting the filetype and fileformat testcases"
#
# The sheet gains a new "StatQuery" column between the existing A (query
# name) and B (dbExcel) columns, and the two query rows are populated
# with the actual Cypher statements that back the "Aligned RNA reads
# file" test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old column B (and C) one slot to the right, but only across
# the two rows that actually hold data -- this keeps the row/col model
# tight instead of inflating the whole 1,048,576-row column.
$ws.Range("B1:B2").Insert(-4161)

# New header + column formatting (mirror column A's width so they render
# as one merged <col> run in the saved file).
$ws.Range("B1").Value = "StatQuery"
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# New query bodies for row 2 (A2 = RNA-reads-only query, B2 = the
# combined file/case/trial count query). Both inherit the existing
# wrap-text style that used to live on A2 alone.
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_type IN [''Aligned RNA reads file'']  RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
$ws.Range("B2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s , f WHERE f.file_type IN [''Aligned DNA reads file'',''Aligned RNA reads file'',''Index file'',''Variants file''] RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trials'
$ws.Range("B2").WrapText = $ws.Range("A2").WrapText

# Row 2 now has to be tall enough to show the wrapped query text.
$ws.Rows("2:2").RowHeight = 101.5

# Land the saved selection/scroll state on A2, with no frozen/scrolled
# top-left cell.
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
